$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "61.693.20"
$ws.Range("E2").Value = "  +1.17%  "
Set-TextValue "D3" "3.409.28"
$ws.Range("E3").Value = "  +0.63%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "577.13"
$ws.Range("E5").Value = "  +0.93%  "
Set-TextValue "D6" "143.67"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.20%  "
Set-TextValue "D9" "7.63"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -0.45%  "
Set-TextValue "D12" "3.990.68"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  -0.58%  "
Set-TextValue "D14" "27.94"
$ws.Range("E14").Value = "  +0.08%  "
Set-TextValue "D15" "3.413.37"
$ws.Range("E15").Value = "  +1.13%  "
Set-TextValue "D16" "0.0000169"
$ws.Range("E16").Value = "  -1.50%  "
Set-TextValue "D17" "61.765.19"
$ws.Range("E17").Value = "  +1.10%  "
Set-TextValue "D18" "6.14"
$ws.Range("E18").Value = "  +0.80%  "
Set-TextValue "D19" "13.74"
$ws.Range("E19").Value = "  +0.68%  "
Set-TextValue "D20" "9.17"
$ws.Range("E20").Value = "  +2.49%  "
Set-TextValue "D21" "387.93"
$ws.Range("E21").Value = "  +1.01%  "
Set-TextValue "D22" "74.28"
$ws.Range("E22").Value = "  -1.07%  "
Set-TextValue "D23" "0.550"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  +0.02%  "
Set-TextValue "D28" "7.39"
$ws.Range("E28").Value = "  +1.62%  "
Set-TextValue "D29" "7.99"
$ws.Range("E29").Value = "  +0.38%  "
Set-TextValue "D30" "2.15"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  +0.00%  "
Set-TextValue "D33" "23.42"
$ws.Range("E33").Value = "  +0.84%  "
Set-TextValue "D34" "6.94"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  +3.49%  "
Set-TextValue "D36" "168.80"
$ws.Range("E36").Value = "  +1.07%  "
Set-TextValue "D37" "3.444.03"
$ws.Range("E37").Value = "  +0.69%  "
Set-TextValue "D38" "1.47"
$ws.Range("E38").Value = "  +0.07%  "
Set-TextValue "D39" "27.99"
$ws.Range("E39").Value = "  +4.08%  "
Set-TextValue "D40" "0.0756"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("E41").Value = "  +0.73%  "
Set-TextValue "D42" "4.44"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("E43").Value = "  +0.50%  "
Set-TextValue "D44" "1.16"
$ws.Range("E44").Value = "  +3.03%  "
Set-TextValue "D45" "2.492.20"
$ws.Range("E45").Value = "  +1.50%  "
Set-TextValue "D46" "22.84"
$ws.Range("E46").Value = "  -0.55%  "
Set-TextValue "D47" "6.63"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  -3.27%  "
$ws.Range("E51").Value = "  -0.39%  "
